# Commit: "added extension delete test"
# Update the "PMTestData" sheet (IP extension test data) to:
#  - extend the existing "create" test's command string (row 13)
#  - add a "delete" test row (row 14)
#  - add a "delete multiple" test row (row 15)
# Also bump the IP addresses used on the "IP" sheet from .113 to .213.

$wb = $excel.ActiveWorkbook
$wsPM = $wb.Worksheets.Item("PMTestData")
$wsIP = $wb.Worksheets.Item("IP")

# --- Row 13: extend the existing command string and grow the row height
#     to fit the now-longer, wrapped text (2 lines -> 4 lines).
$cmdCreate = "number_initiate -number 77777..77779 -numbertype ex,77777-77779,77777,1,FirstName,LastName,Mitel 6869i,ip_extension -e -d 77777,extension -e -d 77777,number_end -number 77777..77779 -numbertype ex"

$wsPM.Range("B13").Value2 = $cmdCreate
$wsPM.Rows.Item(13).RowHeight = 58

# --- "IP" sheet: bump the host IP address used by the test data.
$wsIP.Range("A1").Value2 = "http://10.211.162.213/mp"
$wsIP.Range("A2").Value2 = "http://10.211.162.213/wbm"

# --- Row 14 (new): single IP-extension delete test.
$wsPM.Range("A14").Value2 = "test_delete_IP_extension"
$wsPM.Range("B14").Value2 = $cmdCreate
$wsPM.Range("C14").Value2 = "Y"
$wsPM.Range("B14").WrapText = $true
$wsPM.Rows.Item(14).RowHeight = 58

# --- Row 15 (new): multiple IP-extension delete test.
$cmdDeleteMultiple = "number_initiate -number 70001..70010 -numbertype ex,extension -i -d 70001..70010 -l 1 --csp 0,ip_extension -i -d 70001..70010,70001-70010,ip_extension -e -d 70001..70010,extension -e -d 70001..70010,number_end -number 70001..70010 -numbertype ex"

$wsPM.Range("A15").Value2 = "test_delete_multiple_IP_extension"
$wsPM.Range("B15").Value2 = $cmdDeleteMultiple
$wsPM.Range("C15").Value2 = "Y"
$wsPM.Range("B15").WrapText = $true
$wsPM.Rows.Item(15).RowHeight = 58

# --- Restore/update the view selections on both sheets (IP first, then
#     PMTestData, so PMTestData ends up the active tab, matching the
#     original workbook's activeTab).
$wsIP.Activate()
$wsIP.Range("G3").Select() | Out-Null

$wsPM.Activate()
$wsPM.Range("F15").Select() | Out-Null
